# This workbook lists observation records (rows 4-12 of the active sheet).
# The underlying source data was re-synced and, as a result, the per-record
# fields (Id, Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor,
# Ost, Nord) were reshuffled across the row positions, while the
# Taxonsorteringsordning (column B) values were refreshed to new figures,
# and a "hack" comment moved from row 12 to row 10 (together with a handful
# of blank helper cells in columns K-N).
#
# Below, each row's final (post-edit) values are written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 4;  A = 111936776; B = 77636; D = "NT"; E = 6425;   F = "Garnlav";     G = "Alectoria sarmentosa";       H = "(Ach.) Ach.";                           Q = 490398; R = 7088445 },
    @{ Row = 5;  A = 111936780; B = 77636; D = "NT"; E = 6425;   F = "Garnlav";     G = "Alectoria sarmentosa";       H = "(Ach.) Ach.";                           Q = 489952; R = 7088557 },
    @{ Row = 6;  A = 111936777; B = 77636; D = "NT"; E = 6425;   F = "Garnlav";     G = "Alectoria sarmentosa";       H = "(Ach.) Ach.";                           Q = 490056; R = 7088709 },
    @{ Row = 7;  A = 111936768; B = 90221; D = "LC"; E = 3298;   F = "Trådticka";   G = "Climacocystis borealis";     H = "(Fr.) Kotl. & Pouzar";                  Q = 490317; R = 7088522 },
    @{ Row = 8;  A = 111936781; B = 89927; D = "LC"; E = 4217;   F = "Blodticka";   G = "Meruliopsis taxicola";       H = "(Pers.:Fr.) Bondartsev";                Q = 490315; R = 7088552 },
    @{ Row = 9;  A = 111936779; B = 77636; D = "NT"; E = 6425;   F = "Garnlav";     G = "Alectoria sarmentosa";       H = "(Ach.) Ach.";                           Q = 490008; R = 7088597 },
    @{ Row = 10; A = 111936774; B = 56446; D = "NT"; E = 100049; F = "Spillkråka";  G = "Dryocopus martius";          H = "(Linnaeus, 1758)";                      Q = 490378; R = 7088551 },
    @{ Row = 11; A = 111936767; B = 90221; D = "LC"; E = 3298;   F = "Trådticka";   G = "Climacocystis borealis";     H = "(Fr.) Kotl. & Pouzar";                  Q = 490377; R = 7088412 },
    @{ Row = 12; A = 111936775; B = 89553; D = "NT"; E = 1204;   F = "Gränsticka";  G = "Phellopilus nigrolimitatus"; H = "(Romell) Niemelä, T.Wagner & M.Fisch."; Q = 490380; R = 7088379 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
}

# The "Publik kommentar" (hack) note travels with the record that moves
# from row 12 into row 10; row 12 no longer carries a comment.
$ws.Range("AC10").Value = "hack"
$ws.Range("AC12").ClearContents()
